$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1): B1:D1 ----
$ws.Range("B1").Value = 4
$ws.Range("C1").Value = 5
$ws.Range("D1").Value = 6

# ---- Row 2 ----
$ws.Range("A2").Value = "carID"
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 6

# ---- Row 3 ----
$ws.Range("A3").Value = "speed2"
$ws.Range("B3").Value = 52.2
$ws.Range("C3").Value = 42.72
$ws.Range("D3").Value = 61.85

# ---- Row 4 label ----
$ws.Range("A4").Value = "asma"

# ---- Row 5 ----
$ws.Range("A5").Value = "ceza_tutar"
$ws.Range("B5").Value = 3136
$ws.Range("C5").Value = 1508.5
$ws.Range("D5").Value = 6440

# ---- Row 6 ----
$ws.Range("A6").Value = "hesaplanan_asma"
$ws.Range("B6").Value = 30
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 50

# ---- Styling: bold font, centered/top aligned, thin box border ----
# applies to B1:D1 (header numbers) and A2:A6 (row labels).
# Do this BEFORE the row-4 text trick below so this format lands on
# cellXfs index 1.
$first = $ws.Range("B1")
$first.Font.Bold = $true
$first.Borders.LineStyle = 1
$first.Borders.Weight = 2
$first.HorizontalAlignment = -4108   # xlCenter
$first.VerticalAlignment = -4160     # xlTop

$first.Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A2:A6").PasteSpecial(-4122)   # xlPasteFormats

# ---- Row 4 numeric-looking values, stored as literal text ----
# Leading apostrophe forces text storage instead of numeric parsing;
# resetting the style back to "Normal" strips the quote-prefix
# formatting mark so the cell keeps the default (unstyled) look.
$ws.Range("B4").Value = "'49.15"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'22.06"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "'76.70"
$ws.Range("D4").Style = "Normal"
